$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents for rows 28-51 per commit "add getting sensor AVG"
$ws.Range("A28").Value = "W311b"
$ws.Range("B28").Value = 45089
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = 0.70833333333333337
$ws.Range("E28").Value = "ICU310"
$ws.Range("F28").Value = "Kelvin Leung"
$ws.Range("A29").Value = "W311-H1"
$ws.Range("B29").Value = 45089
$ws.Range("C29").Value = 0.54166666666666663
$ws.Range("D29").Value = 0.70833333333333337
$ws.Range("E29").Value = "ICU310"
$ws.Range("F29").Value = "Eric Tsang"
$ws.Range("A30").Value = "W311-H2"
$ws.Range("B30").Value = 45089
$ws.Range("C30").Value = 0.54166666666666663
$ws.Range("D30").Value = 0.70833333333333337
$ws.Range("E30").Value = "ICU310"
$ws.Range("F30").Value = "Eric Tsang"
$ws.Range("A31").Value = "W311-H3"
$ws.Range("B31").Value = 45089
$ws.Range("C31").Value = 0.54166666666666663
$ws.Range("D31").Value = 0.70833333333333337
$ws.Range("E31").Value = "ICU310"
$ws.Range("F31").Value = "Eric Tsang"
$ws.Range("A32").Value = "W311a"
$ws.Range("B32").Value = 45090
$ws.Range("C32").Value = 0.33333333333333331
$ws.Range("D32").Value = 0.5
$ws.Range("E32").Value = "ICU310"
$ws.Range("F32").Value = "Eric Tsang"
$ws.Range("A33").Value = "W311b"
$ws.Range("B33").Value = 45090
$ws.Range("C33").Value = 0.33333333333333331
$ws.Range("D33").Value = 0.5
$ws.Range("E33").Value = "ICU310"
$ws.Range("F33").Value = "Kelvin Leung"
$ws.Range("A34").Value = "W311b"
$ws.Range("B34").Value = 45090
$ws.Range("C34").Value = 0.5
$ws.Range("D34").Value = 0.70833333333333337
$ws.Range("E34").Value = "ICU310"
$ws.Range("F34").Value = "Kelvin Leung"
$ws.Range("A35").Value = "W311-H1"
$ws.Range("B35").Value = 45090
$ws.Range("C35").Value = 0.33333333333333331
$ws.Range("D35").Value = 0.5
$ws.Range("E35").Value = "IC382"
$ws.Range("F35").Value = "Peter Lam"
$ws.Range("A36").Value = "W311-H1"
$ws.Range("B36").Value = 45090
$ws.Range("C36").Value = 0.54166666666666663
$ws.Range("D36").Value = 0.70833333333333337
$ws.Range("E36").Value = "IC382"
$ws.Range("F36").Value = "Peter Lam"
$ws.Range("A37").Value = "W311-H2"
$ws.Range("B37").Value = 45090
$ws.Range("C37").Value = 0.33333333333333331
$ws.Range("D37").Value = 0.5
$ws.Range("E37").Value = "IC382"
$ws.Range("F37").Value = "Peter Lam"
$ws.Range("A38").Value = "W311-H2"
$ws.Range("B38").Value = 45090
$ws.Range("C38").Value = 0.54166666666666663
$ws.Range("D38").Value = 0.70833333333333337
$ws.Range("E38").Value = "IC382"
$ws.Range("F38").Value = "Peter Lam"
$ws.Range("A39").Value = "W311-H3"
$ws.Range("B39").Value = 45090
$ws.Range("C39").Value = 0.33333333333333331
$ws.Range("D39").Value = 0.5
$ws.Range("E39").Value = "IC382"
$ws.Range("F39").Value = "Peter Lam"
$ws.Range("A40").Value = "W311-H3"
$ws.Range("B40").Value = 45090
$ws.Range("C40").Value = 0.54166666666666663
$ws.Range("D40").Value = 0.70833333333333337
$ws.Range("E40").Value = "IC382"
$ws.Range("F40").Value = "Peter Lam"
$ws.Range("A41").Value = "W311b"
$ws.Range("B41").Value = 45091
$ws.Range("C41").Value = 0.33333333333333331
$ws.Range("D41").Value = 0.5
$ws.Range("E41").Value = "ICU310"
$ws.Range("F41").Value = "Kelvin Leung"
$ws.Range("A42").Value = "W311d-Z1"
$ws.Range("B42").Value = 45091
$ws.Range("C42").Value = 0.33333333333333331
$ws.Range("D42").Value = 0.5
$ws.Range("E42").Value = "ICU310"
$ws.Range("F42").Value = "Kelvin Leung"
$ws.Range("A43").Value = "W311d-Z1"
$ws.Range("B43").Value = 45091
$ws.Range("C43").Value = 0.5
$ws.Range("D43").Value = 0.70833333333333337
$ws.Range("E43").Value = "ICU310"
$ws.Range("F43").Value = "Kelvin Leung"
$ws.Range("A44").Value = "W311d-Z2"
$ws.Range("B44").Value = 45091
$ws.Range("C44").Value = 0.33333333333333331
$ws.Range("D44").Value = 0.5
$ws.Range("E44").Value = "ICU310"
$ws.Range("F44").Value = "KS Chan"
$ws.Range("A45").Value = "W311d-Z2"
$ws.Range("B45").Value = 45091
$ws.Range("C45").Value = 0.5
$ws.Range("D45").Value = 0.70833333333333337
$ws.Range("E45").Value = "ICU310"
$ws.Range("F45").Value = "KS Chan"
$ws.Range("A46").Value = "W311-H1"
$ws.Range("B46").Value = 45091
$ws.Range("C46").Value = 0.33333333333333331
$ws.Range("D46").Value = 0.5
$ws.Range("E46").Value = "IC382"
$ws.Range("F46").Value = "Peter Lam"
$ws.Range("A47").Value = "W311-H1"
$ws.Range("B47").Value = 45091
$ws.Range("C47").Value = 0.54166666666666663
$ws.Range("D47").Value = 0.70833333333333337
$ws.Range("E47").Value = "IC382"
$ws.Range("F47").Value = "Peter Lam"
$ws.Range("A48").Value = "W311-H2"
$ws.Range("B48").Value = 45091
$ws.Range("C48").Value = 0.33333333333333331
$ws.Range("D48").Value = 0.5
$ws.Range("E48").Value = "IC382"
$ws.Range("F48").Value = "Peter Lam"
$ws.Range("A49").Value = "W311-H2"
$ws.Range("B49").Value = 45091
$ws.Range("C49").Value = 0.54166666666666663
$ws.Range("D49").Value = 0.70833333333333337
$ws.Range("E49").Value = "IC382"
$ws.Range("F49").Value = "Peter Lam"
$ws.Range("A50").Value = "W311-H3"
$ws.Range("B50").Value = 45091
$ws.Range("C50").Value = 0.33333333333333331
$ws.Range("D50").Value = 0.5
$ws.Range("E50").Value = "IC382"
$ws.Range("F50").Value = "Peter Lam"
$ws.Range("A51").Value = "W311-H3"
$ws.Range("B51").Value = 45091
$ws.Range("C51").Value = 0.54166666666666663
$ws.Range("D51").Value = 0.70833333333333337
$ws.Range("E51").Value = "IC382"
$ws.Range("F51").Value = "Peter Lam"

# Remove now-obsolete trailing rows (52-53)
$ws.Rows("52:53").Delete()

# Refresh AutoFilter range to the new data extent
$ws.AutoFilterMode = $false
$ws.Range("A2:F51").AutoFilter()

# Re-apply the sort so sortState reflects the new extent
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B3:B51"))
$ws.Sort.SetRange($ws.Range("A2:F51"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Keep the _FilterDatabase defined name in sync with the new range
$name = $wb.Names.Item(1)
$name.RefersTo = "='Event for Group A'!`$A`$2:`$F`$51"

# Update the saved selection to match the author's final cursor position
$ws.Range("H29").Select()
